# Weekly update: insert a new daily price record at row 42 (Vega Monumental
# Concepción - Mango), pushing the existing rows 42-71 down to 43-72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 42; everything below shifts
# down by one (old row 42 becomes row 43, ..., old row 71 becomes row 72).
$ws.Rows.Item(42).EntireRow.Insert()

# Populate the newly inserted row 42 with the latest weekly observation.
$ws.Range("A42").Value = 11
$ws.Range("B42").Value = "Vega Monumental Concepción"
$ws.Range("C42").Value = "Bíobío"
$ws.Range("D42").Value = 44483
$ws.Range("E42").Value = 8
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100108
$ws.Range("H42").Value = "Tropicales y subtropicales"
$ws.Range("I42").Value = 100108002
$ws.Range("J42").Value = "Mango"
$ws.Range("K42").Value = "Sin especificar"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = 7000
$ws.Range("O42").Value = 7500
$ws.Range("P42").Value = 7250
$ws.Range("Q42").Value = "$/bandeja 4 kilos"
$ws.Range("R42").Value = "Perú"
$ws.Range("S42").Value = 1812
$ws.Range("T42").Value = 4
